$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row for the "NameDrei" group-of-three IBAN payment that used
# to be recorded in row 4; this pushes it down to row 5 and frees up row 4
# for the corrected, now-complete Maringer "Buchung".
$ws.Rows.Item(5).Insert()

# Row 3 (Corelie Scholz): booking date corrected 21/09/2016 -> 23/09/2016
# and amount corrected 130 -> 120.
$ws.Range("A3").Value = 42636
$ws.Range("E3").Value = 120

# Row 4 (Johannes MARINGER, full Buchung incl. IBAN-matched amount):
# booking date corrected 03/10/2016 -> 21/09/2016, amount now 130.
$ws.Range("A4").Value = 42634
$ws.Range("E4").Value = 130

# Row 5 (new): "NameDrei" group-of-three payment, 03/10/2016, amount 70.
$ws.Range("A5").Value = 42646
$ws.Range("C5").Value = 10
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Value = "Gruppenbeitrag 16/17 NameDrei VornameDrei1 VornameDrei2 VornameDrei3"
$ws.Range("E5").Value = 70
